$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.888.74"
$ws.Range("E2").Value = "  +1.92%  "

$ws.Range("D3").Value = "3.056.21"
$ws.Range("E3").Value = "  +0.85%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.76%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  +0.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.29"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.06%  "

$ws.Range("E10").Value = "  -0.70%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.374"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.00%  "

$ws.Range("D12").Value = "3.578.37"
$ws.Range("E12").Value = "  +0.78%  "

$ws.Range("E13").Value = "  +2.97%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.26%  "

$ws.Range("E15").Value = "  +0.10%  "

$ws.Range("D16").Value = "57.886.22"
$ws.Range("E16").Value = "  +1.91%  "

$ws.Range("D17").Value = "3.052.90"
$ws.Range("E17").Value = "  +0.60%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.11%  "

$ws.Range("E20").Value = "  -0.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "331.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.91%  "

$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("E23").Value = "  -0.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.03%  "

$ws.Range("E25").Value = "  +2.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").Value = "0.0₃0898"
$ws.Range("E27").Value = "  -3.64%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.22%  "

$ws.Range("E30").Value = "  +1.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.72%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.49%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.71%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.75%  "

$ws.Range("E37").Value = "  +2.82%  "

$ws.Range("E38").Value = "  +2.29%  "

$ws.Range("D39").Value = "3.094.68"
$ws.Range("E39").Value = "  +0.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.77%  "

$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("E42").Value = "  -0.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.653"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.29%  "

$ws.Range("D44").Value = "2.270.53"
$ws.Range("E44").Value = "  +2.62%  "

$ws.Range("E45").Value = "  +5.36%  "

$ws.Range("E46").Value = "  +1.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.69%  "

$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.937"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.07%  "

$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.735"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.09%  "

$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0877"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.40%  "
